$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (emails are no longer mailto-linked)
$ws.Hyperlinks.Delete()

# Update names (column B) and emails (column C) for each student row
$ws.Range("B2").Value = "Darshan"
$ws.Range("C2").Value = "pasne.d@husky.neu.edu"

$ws.Range("B3").Value = "Saman"
$ws.Range("C3").Value = "sood.s@husky.neu.edu"

$ws.Range("B4").Value = "Shail"
$ws.Range("C4").Value = "shail@ccs.neu.edu"

$ws.Range("B5").Value = "Vaibhav"
$ws.Range("C5").Value = "dave.v@husky.neu.edu"

$ws.Range("B6").Value = "John"
$ws.Range("C6").Value = "snow.j@husky.neu.edu"

$ws.Range("B7").Value = "Danny"
$ws.Range("C7").Value = "danny.d@husky.neu.edu"

$ws.Range("B8").Value = "Erica"
$ws.Range("C8").Value = "sniper.e@husky.neu.edu"

$ws.Range("B9").Value = "Flurry"
$ws.Range("C9").Value = "majin.f@husky.neu.edu"

$ws.Range("B10").Value = "Gara"
$ws.Range("C10").Value = "hawking.g@husky.neu.edu"

$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"

$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"

# Move selection to reflect the author's final cursor position
$ws.Range("C19").Select() | Out-Null
